# Generate Report for Handoff
# - Flip status cells from "In Translation" to "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the status columns to fit the new, longer status text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-10-20 00:54:52"
$wsZhCn.Range("H2").Value     = "2016-10-20 00:54:41"
$wsDeDe.Range("H2").Value     = "2016-10-20 00:54:52"

# --- Column widths (status columns widened for the longer text) ---
# The stored OOXML "width" attribute = ColumnWidth + 5/6, so subtract 5/6
# from the desired stored width before assigning it through COM.
$newStatusColWidth = 17.2159881591797 - (5/6)

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newStatusColWidth
